$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "56.872.16"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.38%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.051.76"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.87%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "518.56"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.49%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.28"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.65%  "

$ws.Range("E7").Value = "  +0.10%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.050.45"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.87%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.462"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.55%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.25"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.47%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.105"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.75%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.400"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.71%  "

$ws.Range("E13").Value = "  +0.93%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.579.09"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.80%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "24.90"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.77%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000158"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.51%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "56.917.42"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.36%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.057.10"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.69%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.81"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.85%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.30"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.01%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.71"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.46%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "345.42"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.68%  "

$ws.Range("E23").Value = "  -0.07%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "68.43"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.13%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.494"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.81%  "

$ws.Range("E26").Value = "  +0.29%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.162"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.21%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0836"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -10.18%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.08"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.18%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.84"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.44%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.74"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.76%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.71"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -11.06%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "157.65"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.53%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.76"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.05%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.10"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -7.21%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.92"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.60%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "25.01"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.10%  "

$ws.Range("E39").Value = "  -2.97%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0649"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.90%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.54"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -6.87%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.97"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.60%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.685"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.48%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.390.67"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.98%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "36.34"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.01%  "

$ws.Range("E46").Value = "  -0.01%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.093.96"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.70%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0258"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.17%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.93"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.81%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.918"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -9.09%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.17"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -7.28%  "
